# Booking_Creation_DataSet.xlsx -- add end-to-end test rows (E2E + UPDATE)
# and validate the previously last "BOOKING_VALID" row's roomid value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 60 (BOOKING_VALID): bump the roomid from 356 to 588 -----------
$ws.Cells.Item(60, 2).Value = 588

# --- 2. Row 61 (new BOOKING_E2E row) ---------------------------------------
# Clone row 60's full formatting (borders/fills/number-formats) down onto
# row 61 first, then overwrite the cells that actually differ.
$ws.Range("A60:K60").Copy($ws.Range("A61"))

$ws.Cells.Item(61, 1).Value = "BOOKING_E2E"
$ws.Cells.Item(61, 2).Value = 4652
$ws.Cells.Item(61, 3).Value = "Saravanan S"
$ws.Cells.Item(61, 4).Value = "Subramaniyan"
$ws.Cells.Item(61, 5).Value = "true"
$ws.Cells.Item(61, 6).Value = 46019
$ws.Cells.Item(61, 7).Value = 46022
$ws.Cells.Item(61, 8).Value = "subbusrvn@gmail.com"
$ws.Cells.Item(61, 9).Value = "919710288425"
$ws.Cells.Item(61, 10).Value = "Schema Validation Check"
$ws.Cells.Item(61, 11).Value = "Booking should be created"

# Column A on this row lost its outline border in the authored workbook
# (the cell style that A61 ended up with has borderId 0) -- match that.
$ws.Cells.Item(61, 1).Borders.LineStyle = -4142

# --- 3. Row 62 (new BOOKING_UPDATE row) ------------------------------------
$ws.Range("A60:K60").Copy($ws.Range("A62"))

$ws.Cells.Item(62, 1).Value = "BOOKING_UPDATE"
$ws.Cells.Item(62, 2).Value = 5875
$ws.Cells.Item(62, 3).Value = "Samyuktha"
$ws.Cells.Item(62, 4).Value = "Saravanan_Update"
$ws.Cells.Item(62, 5).Value = "false"
$ws.Cells.Item(62, 6).Value = 46019
$ws.Cells.Item(62, 7).Value = 46022
$ws.Cells.Item(62, 8).Value = "samsaravanan@gmail.com"
$ws.Cells.Item(62, 9).Value = "919710288178"
$ws.Cells.Item(62, 10).Value = ""
$ws.Cells.Item(62, 11).Value = "Booking should be updated"
